$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33
$ws.Cells.Item($row, 1).Value = "Golang Architect / Principal Backend Architect"
$ws.Cells.Item($row, 2).Value = "https://www.dice.com/job-detail/679f7586-7818-4fbb-9e79-a94c62a72734"
$ws.Cells.Item($row, 3).Value = "Georgia"
$ws.Cells.Item($row, 4).Value = "Third Party"
$ws.Cells.Item($row, 5).Value = "`$80 - `$85"
$ws.Cells.Item($row, 6).Value = "Source Mantra Inc"
